$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B14").Value = "U2FsdGVkX19pVFiEFGr+okGR3DwaGQRgeHdGXpjw4VUK9xa0KlweeGwh3FMR3QoANLgsNWRuwDHrRafszy90digd2jmAQyBN3iwLMoIkK/lraSwhNIMD/Fox+vsCEHHcDhwRD+XEFreFcoVxCsRav3DDT3GwCkwpu3HMku2sWJ+VZF4/075wIvi+6DN8j46gTe0uXl59ZBOwleRKiDjon5mMq1KTl3vMn1rKvKzyLflmyrXLM4ZCvcvI3vxdQ6bd6D/+3mbvJt6KYZubNl//Qi2YZ/nBsntb7U+PZKlpwDO8lpNO3O4d+UTvOQAwoao3YRrKY7AvJ+SyNTTYUbTiI2OGIoqpQMNirm36S/effJddVLJn9KZInBPXlAYJoP0Uv/zhNklnQdYl88DMxfMRHPbUA4fobpuSwz5XYPOlL6MxhpA9+CQppZ2CYcuPi0K2PZLyerGnyI48wbr+XR3B/h7bi3gXGS3O1MQZGX5Vw6kiZf3iPtnzmBO8RHoeKEaq8+MvM/V/VjV0sOrZLWOJ02KFi6Nx3Kl1rVruaNuCPLOBdolSjKvCzILvR1kx5W8FvrEueogX9TcTNYMQ3Xm43lG1s1/1wf3aoOI/tIA+FQSrPqJ4N2wk5ROCoPplc/KcDVaUMKE1ELkMWW2B9b1Pi7xa/saPjYStUQoqGM31A628kaISNc0mEtg7LPxZ144nsSjfYMpuWWIXyEGe4xbg30JBo/npQZOvmDW2eDnf0xACUBr8x7YJtso8sgbOm+6RQtYTNUXO51sWYHL47vmXMEB1KN9AvG7FsDrjyAIDyjjhX7nVjKcr7uJDTjLA8Ka78zX3DZAUPT3Lt3KKO/ZnfGziSglchoIMQo9rm1IB7u4v8bnpK7wTxyTufjHj/WEkLLChtw82rLDqSFvl5PmiY1pULRwNU3BwtEpkLLuvBjpOH7oDKtY8i+SDjJQWKK2bMJd581Ym6+TYEVVOu32mokZqXuzBAvxVgs7pT3Zr/NfJbfJHrMzAfXZXAbbo6jbLROlSBZR+M51fY+qySslWKL+aIDurEiH5wYHkWzKs+hdiixcNhToD6atdnMwsjbVE"
